$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($RowNum, $D, $J, $K, $L, $M, $P) {
    $ws.Cells.Item($RowNum, 1).Value = 1
    $ws.Cells.Item($RowNum, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($RowNum, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($RowNum, 4).Value = $D
    $ws.Cells.Item($RowNum, 5).Value = 15
    $ws.Cells.Item($RowNum, 6).Value = 100112038
    $ws.Cells.Item($RowNum, 7).Value = "Cebollín baby"
    $ws.Cells.Item($RowNum, 8).Value = "Sin especificar"
    $ws.Cells.Item($RowNum, 9).Value = "Primera"
    $ws.Cells.Item($RowNum, 10).Value = $J
    $ws.Cells.Item($RowNum, 11).Value = $K
    $ws.Cells.Item($RowNum, 12).Value = $L
    $ws.Cells.Item($RowNum, 13).Value = $M
    $ws.Cells.Item($RowNum, 14).Value = "$/paquete 1,5 a 2 kilos"
    $ws.Cells.Item($RowNum, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($RowNum, 16).Value = $P
    $ws.Cells.Item($RowNum, 17).Value = 2
    $ws.Cells.Item($RowNum, 18).Value = "Hortaliza"
}

# 1) Insert a brand-new weekly record above the former row 17 (new week's data),
#    pushing all the existing rows down by one.
$ws.Rows.Item(17).Insert()
Set-DataRow 17 44428 270 3500 3800 3650 1825

# 2) Insert two more new weekly records above the row that now holds the old
#    row-27 data (date 44405), pushing it (and everything after) down by two.
$ws.Rows.Item(28).Insert()
Set-DataRow 28 44435 500 1800 2000 1930 965

$ws.Rows.Item(29).Insert()
Set-DataRow 29 44431 300 1900 2000 1950 975

# 3) Insert one more new weekly record above the row that now holds the old
#    row-31 data (date 44319), pushing it (and everything after) down by one.
$ws.Rows.Item(34).Insert()
Set-DataRow 34 44433 200 1800 2000 1900 950
